# Update the "想去人数" (want-to-go count) and "最低票价" (lowest price)
# figures on the "展览" and "全部类型" sheets to the newly scraped values.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 6728
    $ws.Range("G2").Value = 58

    $ws.Range("F3").Value = 46

    $ws.Range("F5").Value = 1052

    $ws.Range("F6").Value = 144
}
